$d = $word.ActiveDocument

# Find the paragraph containing the sentence that needs to be removed, and
# delete the whole paragraph (including its trailing paragraph mark) so that
# the surrounding paragraphs merge back together seamlessly.
$target = "כאשר הצמתים שמתאימים למשתנים הראשיים הם נקודות ההתחלה של שוויונות בין צמתים."

foreach ($p in $d.Paragraphs) {
    $r = $p.Range
    if ($r.Text -like "*$target*") {
        $r.Delete()
        break
    }
}
